# =============================================================================
# edit.ps1
#
# Implements the diff:
#  1. Replaces the short "Ponizej znajduje sie przykladowy plik..." paragraph
#     (after "Kontener 5 wykorzystuje...") with a full Docker Compose
#     write-up: an introduction, the docker-compose.yml listing (many short
#     single-line paragraphs), an explanation, run instructions and a test
#     description. The _GoBack bookmark that used to sit at the end of the
#     "Kontener 5..." paragraph is removed (it moves further down).
#  2. Removes the (now redundant) empty bold/32pt paragraph before the
#     "Koncepcja i architektura rozwiazania" heading, and turns the empty
#     bold/32pt paragraph that used to follow the heading into a normal,
#     unformatted paragraph containing an introductory sentence.
#  3. Splits the "Wymagania funkcjonalne i niefunkcjonalne" heading
#     paragraph into two separate bold headings ("Wymagania funkcjonalne:"
#     and "Wymagania niefunkcjonalne:") separated by blank spacer
#     paragraphs; the _GoBack bookmark reappears at the start of the new
#     "Wymagania niefunkcjonalne:" heading paragraph.
# =============================================================================

$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

function Get-ParaContaining {
    param($doc, [string]$needle)
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    throw "Paragraph containing '$needle' not found"
}

function New-PkgXml {
    param([string]$bodyFragment)
    $header = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $footer = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $header + $bodyFragment + $footer
}

function Insert-XmlAtCollapsedRange {
    param($rng, [string]$bodyFragment)
    $pkg = New-PkgXml $bodyFragment
    $rng.InsertXML($pkg)
}

# -----------------------------------------------------------------------
# STEP 1 (do this part of the document first - it sits further down, so
# editing it first means paragraph indices used later for the earlier
# part of the document, below, are unaffected):
#
# Split "Wymagania funkcjonalne i niefunkcjonalne" into two headings.
# -----------------------------------------------------------------------

# The _GoBack bookmark currently sits inside the "Kontener 5..." paragraph;
# it will be recreated (at its new location) as part of the XML inserted
# below, so drop the old one now.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$pWymagania = Get-ParaContaining $d "Wymagania funkcjonalne i niefunkcjonalne"
$pWymagania.Range.Text = "Wymagania funkcjonalne:"

$insAfterWymagania = $d.Range($pWymagania.Range.End - 1, $pWymagania.Range.End - 1)
$block3 = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve">Wymagania </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>niefunkcjonalne</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
</w:p>
'@
Insert-XmlAtCollapsedRange $insAfterWymagania $block3

# -----------------------------------------------------------------------
# STEP 2: Insert the Docker Compose write-up after "Kontener 5..." and
# remove the old one-line "Ponizej znajduje sie..." paragraph.
# -----------------------------------------------------------------------

$pKontener5 = Get-ParaContaining $d "Kontener 5 wykorzystuje oficjalny obraz serwera bazy danych redis"

$insAfterKontener5 = $d.Range($pKontener5.Range.End - 1, $pKontener5.Range.End - 1)
$block1 = @'
<w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">By wygodnie uruchomić jednocześnie wszystkie kontenery wykorzystaliśmy narzędzie DockerCompose. Zaczęliśmy od stworzenia następującego pliku konfiguracyjnego </w:t>
      </w:r>
      <w:r>
        <w:t>docker-compose.yml</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>nginx:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    build: ./nginx</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    ports:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">        - "80:80"</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>node1:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    build: ./node</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    volumes:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">     - .:/code</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    depends_on:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">     - redis</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>node2:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    build: ./node</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    volumes:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">     - .:/code</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    depends_on:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">     - redis</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>node3:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    build: ./node</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    volumes:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">     - .:/code</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    depends_on:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">     - redis</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>redis:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">    image: redis</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Konfigurowany jest w nim po kolei każdy kontener. Konfiguracja kontenera zawiera informację gdzie może być znaleziony plik Dockerfile potrzebny do zbudowania potrzebnego obrazu, jakie porty mają być</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> zmapowane z portami fizycznego urządzenia na którym uruchamiane są kontenery. Przy kontenerach odpowiedzialnych za serwery http (node1, node2, node3) wprowadzony jest dodatkowy parametr depends_on, wpływający na to że zarządca program Docker Compose najpierw uruchomi kontener obsługujący serwer bazy danych redis a dopiero w następnej kolejności zależne od niego kontenery node1, node2, node3.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Po przygotowaniu środowiska możemy je uruchomić w trybie interaktywnym, poleceniem:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>docker-compose up</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Gdy wszystkie kontenery zostaną uruchomione, możemy przetestować aplikację, wysyłając zapytanie http na port 80. W odpowiedzi dostaniem zawsze identyfikator następnego węzła z pośród trzech </w:t>
      </w:r>
      <w:r>
        <w:t>odpowiedzialnych za obsługę zapytań http oraz kolejny numer licznika odwiedzin. Pokazuje to prawidłowe działanie komunikujących się ze sobą kontenerów oraz równoważenie obciążenia poprzez load balancer.  Widać również że kontenery korzystają z jednej tej samej bazy danych, ponieważ licznik odwiedzin jest odpowiednio zwiększany po każdym zapytaniu.</w:t>
      </w:r>
    </w:p>
    <w:p/>
'@
Insert-XmlAtCollapsedRange $insAfterKontener5 $block1

$pPoni = Get-ParaContaining $d "Ponizej znajduje sie przykladowy plik konfiguracyjny Docker Compose"
$pPoni.Range.Delete()

# -----------------------------------------------------------------------
# STEP 3: Remove the empty heading-styled paragraph before "Koncepcja i
# architektura rozwiazania", and turn the empty heading-styled paragraph
# after it into a plain paragraph with the new introductory sentence.
# -----------------------------------------------------------------------

$pKoncepcja = Get-ParaContaining $d "Koncepcja i architektura rozwiazania"
$pBefore = $pKoncepcja.Previous()
$pBefore.Range.Delete()

$pKoncepcja = Get-ParaContaining $d "Koncepcja i architektura rozwiazania"
$pAfter = $pKoncepcja.Next()
$pAfter.Range.Delete()

$insAfterKoncepcja = $d.Range($pKoncepcja.Range.End - 1, $pKoncepcja.Range.End - 1)
$naszaGrupa = @'
Nasza grupa zdecydowała się na stworzenie aplikacji pozwalającej jej użytkownikom wymieniać się plikami graficznymi 
'@
$naszaGrupaFrag = '<w:p><w:r><w:t xml:space="preserve">' + $naszaGrupa + '</w:t></w:r></w:p>'
Insert-XmlAtCollapsedRange $insAfterKoncepcja $naszaGrupaFrag
